# Slide 2 ("The CNN program"), Content Placeholder 2: the sentence lists the
# libraries used by the CNN project. The author added "Tensorflow, " in front
# of the existing "Keras" mention, e.g.:
#   ... and uses the Keras and Theano libraries.
# becomes
#   ... and uses the Tensorflow, Keras and Theano libraries.
#
# In the canonical OOXML this shows up as two brand-new <a:r> runs inserted
# right before the existing "Keras" run:
#   <a:r><a:rPr lang="en-US" dirty="0" err="1"/><a:t>Tensorflow</a:t></a:r>
#   <a:r><a:rPr lang="en-US"/><a:t>, </a:t></a:r>
#   <a:r><a:rPr lang="en-US" dirty="0" err="1"/><a:t>Keras</a:t></a:r>   <- unchanged original run

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Locate the existing "Keras" run's text and grab it as a sub-range.
$full = $tr.Text
$idx  = $full.IndexOf("Keras")
$kerasRange = $tr.Characters($idx + 1, 5)

# Insert "Tensorflow, " right before it. PowerPoint merges the freshly typed
# text into the "Keras" run because it matches its character formatting
# (same err="1" misspelling flag), giving one run: "Tensorflow, Keras".
$kerasRange.InsertBefore("Tensorflow, ") | Out-Null

# Now split "Tensorflow, Keras" back into three runs by touching only the
# middle ", " slice. Re-assigning a (no-op) character-formatting property on
# that sub-range forces the engine to materialize it as its own run, without
# disturbing the "Tensorflow" and "Keras" runs on either side.
$full2 = $tr.Text
$idx2  = $full2.IndexOf("Tensorflow, Keras")
$commaRange = $tr.Characters($idx2 + 1 + "Tensorflow".Length, 2)
$commaRange.Font.Kerning = $commaRange.Font.Kerning

Write-Output $tr.Text
